$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 29 new rows (40-68) of MarketTitle / MarketId pairs sourced from TSETMC,
# using the exact (unnormalized) Arabic-Yeh strings from the source, distinct from
# the existing normalized Farsi-Yeh entries already present in rows 2-39.
$ws.Cells.Item(40, 1).Value = 'بازار ابزارهاي نوين مالي فرابورس'
$ws.Cells.Item(40, 2).Value = 'IFB.NFinTools'
$ws.Cells.Item(41, 1).Value = 'بازار پايه زرد فرابورس'
$ws.Cells.Item(41, 2).Value = 'IFB.Paye.Yellow'
$ws.Cells.Item(42, 1).Value = 'بازار پايه زرد فرابورس قانون احکام دائمی برنامه های توسعه کشور'
$ws.Cells.Item(42, 2).Value = 'IFB.Paye.Yellow'
$ws.Cells.Item(43, 1).Value = 'بازار اول (تابلوي اصلي) بورس'
$ws.Cells.Item(43, 2).Value = 'TSE.M1.Main'
$ws.Cells.Item(44, 1).Value = 'بازار اول (تابلوي فرعي) بورس'
$ws.Cells.Item(44, 2).Value = 'TSE.M1.Subsidiary'
$ws.Cells.Item(45, 1).Value = 'بازار پايه نارنجي فرابورس'
$ws.Cells.Item(45, 2).Value = 'IFB.Paye.Orange'
$ws.Cells.Item(46, 1).Value = 'بازار ابزارهاي مشتقه فرابورس'
$ws.Cells.Item(46, 2).Value = 'IFB.Derivatives'
$ws.Cells.Item(47, 1).Value = 'بازار اوراق بدهي'
$ws.Cells.Item(47, 2).Value = 'TSE.Bonds'
$ws.Cells.Item(48, 1).Value = 'بازار پايه زرد فرابورس لغو پذیرش شده'
$ws.Cells.Item(48, 2).Value = 'IFB.Paye.Yellow'
$ws.Cells.Item(49, 1).Value = 'بازار پايه نارنجي فرابورس لغو پذیرش شده'
$ws.Cells.Item(49, 2).Value = 'IFB.Paye.Orange'
$ws.Cells.Item(50, 1).Value = 'بازار پايه قرمز فرابورس قانون احکام دائمی برنامه های توسعه کشور'
$ws.Cells.Item(50, 2).Value = 'IFB.Paye.Red'
$ws.Cells.Item(51, 1).Value = 'بازار سوم فرابورس قانون احکام دائمی برنامه های توسعه کشور'
$ws.Cells.Item(51, 2).Value = 'IFB.M3'
$ws.Cells.Item(52, 1).Value = 'بازار پايه قرمز فرابورس قانون احکام دائمی برنامه های توسعه کشور/ اظهارنظر مردود'
$ws.Cells.Item(52, 2).Value = 'IFB.Paye.Red'
$ws.Cells.Item(53, 1).Value = 'بازار عادي آتي'
$ws.Cells.Item(53, 2).Value = 'TSE.Derivatives'
$ws.Cells.Item(54, 1).Value = 'بازار پايه نارنجي فرابورس قانون احکام دائمی برنامه های توسعه کشور/عدم اظهارنظر'
$ws.Cells.Item(54, 2).Value = 'IFB.Paye.Orange'
$ws.Cells.Item(55, 1).Value = 'بازار ابزارهاي مشتقه'
$ws.Cells.Item(55, 2).Value = 'IFB.Derivatives'
$ws.Cells.Item(56, 1).Value = 'بازار پايه نارنجي فرابورس لغو پذیرش شده/ اظهارنظر مردود'
$ws.Cells.Item(56, 2).Value = 'IFB.Paye.Orange'
$ws.Cells.Item(57, 1).Value = 'بازار پايه قرمز فرابورس لغو پذیرش شده/اظهارنظر مردود'
$ws.Cells.Item(57, 2).Value = 'IFB.Paye.Red'
$ws.Cells.Item(58, 1).Value = 'بازار پايه نارنجي فرابورس لغو پذیرش شده/عدم اظهارنظر'
$ws.Cells.Item(58, 2).Value = 'IFB.Paye.Orange'
$ws.Cells.Item(59, 1).Value = 'بازار پايه  فرابورس'
$ws.Cells.Item(59, 2).Value = 'IFB.Paye'
$ws.Cells.Item(60, 1).Value = 'شرکتهاي کوچک و متوسط فرابورس'
$ws.Cells.Item(60, 2).Value = 'IFB.SME'
$ws.Cells.Item(61, 1).Value = 'بازار پايه قرمز فرابورس'
$ws.Cells.Item(61, 2).Value = 'IFB.Paye.Red'
$ws.Cells.Item(62, 1).Value = 'بازار پايه قرمز فرابورس لغو پذیرش شده/عدم اظهارنظر/ورشکستگی'
$ws.Cells.Item(62, 2).Value = 'IFB.Paye.Red'
$ws.Cells.Item(63, 1).Value = 'بازار دوم فرابورس قانون احکام دائمی برنامه های توسعه کشور'
$ws.Cells.Item(63, 2).Value = 'IFB.M2'
$ws.Cells.Item(64, 1).Value = 'بازار پايه نارنجي فرابورس قانون احکام دائمی برنامه های توسعه کشور'
$ws.Cells.Item(64, 2).Value = 'IFB.Paye.Orange'
$ws.Cells.Item(65, 1).Value = 'بازار پايه  فرابورس قانون احکام دائمی برنامه های توسعه کشور'
$ws.Cells.Item(65, 2).Value = 'IFB.Paye'
$ws.Cells.Item(66, 1).Value = 'بازار پايه نارنجي فرابورس لغو پذیرش شده/ عدم اظهارنظر'
$ws.Cells.Item(66, 2).Value = 'IFB.Paye.Orange'
$ws.Cells.Item(67, 1).Value = 'بازار پايه قرمز فرابورس لغو پذیرش شده/ عدم اظهارنظر/ انحلال'
$ws.Cells.Item(67, 2).Value = 'IFB.Paye.Red'
$ws.Cells.Item(68, 1).Value = 'بازار پايه قرمز فرابورس لغو پذیرش شده'
$ws.Cells.Item(68, 2).Value = 'IFB.Paye.Red'
